# Apply "Correction" edit to FallDDDQN.xlsx
# - Multiply columns I and J (rows 2-22) by 10
# - Change the active selection to P10
# - Auto-fit / set width of columns L:M (bestFit) to 10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the Avg_Throughput_Before / Avg_Throughput_After values (columns I & J) ---
# The source data had values scaled down by a factor of 10; multiply them back.
for ($r = 2; $r -le 22; $r++) {
    $iCell = $ws.Cells.Item($r, 9)
    $jCell = $ws.Cells.Item($r, 10)
    $iCell.Value = $iCell.Value() * 10
    $jCell.Value = $jCell.Value() * 10
}

# --- Columns L:M get a best-fit width of 10 (chars) ---
# (ColumnWidth 9.14 -> stored sheet width of 10, matching Excel's AutoFit result
# for these columns; the COM ColumnWidth setter applies Excel's usual internal
# padding when it serializes the column's stored "width".)
$ws.Range("L1:M1").ColumnWidth = 9.14

# --- Update the active selection to P10 ---
$ws.Range("P10").Select()
